$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: "Database(SQL,MongoDB,SQLite)" -> B6 gets the fulfillment note
$ws.Range("B6").Value = "To be fulfilled by Thursday September 30th 2021"

# Row 7: "JS library not covered in class" -> B7 changes from "Google Charts ?" to "FullChart"
$ws.Range("B7").Value = "FullChart"
